$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44552
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 15500
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15750
$ws.Range("Q2").Value = "$/caja 15 kilos"
$ws.Range("S2").Value = 1050
$ws.Range("T2").Value = 15

# Row 3
$ws.Range("D3").Value = 44165
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("Q3").Value = "$/caja 15 kilos granel"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1100
$ws.Range("T3").Value = 15

# Row 4
$ws.Range("D4").Value = 44187
$ws.Range("K4").Value = "Dina"
$ws.Range("M4").Value = 55
$ws.Range("N4").Value = 15000
$ws.Range("P4").Value = 15455
$ws.Range("Q4").Value = "$/caja 15 kilos granel"
$ws.Range("S4").Value = 1030

# Row 5
$ws.Range("D5").Value = 44189
$ws.Range("K5").Value = "Dina"
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 16562
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("S5").Value = 920
$ws.Range("T5").Value = 18

# Row 7
$ws.Range("D7").Value = 44176
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 17400
$ws.Range("S7").Value = 967

# Row 8
$ws.Range("D8").Value = 44551
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 15500
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15750
$ws.Range("Q8").Value = "$/caja 15 kilos"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1050
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44168
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 16500
$ws.Range("Q9").Value = "$/caja 16 kilos granel"
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 1031
$ws.Range("T9").Value = 16

# Row 10
$ws.Range("D10").Value = 44174
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 75
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 9467
$ws.Range("Q10").Value = "$/caja 10 kilos"
$ws.Range("S10").Value = 947
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44537
$ws.Range("L11").Value = "Primera"
$ws.Range("N11").Value = 21000
$ws.Range("O11").Value = 21500
$ws.Range("P11").Value = 21250
$ws.Range("Q11").Value = "$/caja 15 kilos"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 1417

# Row 12
$ws.Range("D12").Value = 44544
$ws.Range("K12").Value = "Castle Brite"
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16500
$ws.Range("Q12").Value = "$/caja 15 kilos"
$ws.Range("S12").Value = 1100

Write-Output "Updated weekly Damasco price rows"
